# Re-process the metadata sheet with the newly curated dimensions.
# (Se procesan de nuevo los datos con las nuevas dimensiones curadas)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: sdmx/iaest annotation row
$ws.Range("C2").Value = "sdmx-dimension:refArea"
$ws.Range("D2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:edad-grandes-grupos"

# Row 3: dim/medida row
$ws.Range("D3").Value = "medida"
$ws.Range("G3").Value = "medida"

# Row 4: URI/xsd type row
$ws.Range("C4").Value = "URI-Comunidad"
$ws.Range("D4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"

# Row 5 (mapping-*.xlsx helper row) is no longer needed - remove it entirely.
$ws.Rows.Item(5).Delete()
